# Auto-generated edit script: update cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.770.12"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "3.049.33"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'559.04"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'142.29"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.047.98"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").Value = "'0.515"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").Value = "'6.31"
$ws.Range("E11").Value = "  -10.36%  "
$ws.Range("D12").Value = "'0.489"
$ws.Range("E12").Value = "  +5.75%  "
$ws.Range("D13").Value = "'0.0000229"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "'35.66"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "3.551.78"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").Value = "63.841.27"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "3.053.77"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").Value = "'6.79"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "'475.82"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "'14.03"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").Value = "'14.69"
$ws.Range("E22").Value = "  +10.45%  "
$ws.Range("D23").Value = "'0.682"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").Value = "'7.53"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").Value = "'82.82"
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'2.79"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").Value = "'8.12"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "'2.03"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'26.22"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").Value = "'2.44"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "'5.76"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "'6.21"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "'54.52"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").Value = "'0.0409"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'445.23"
$ws.Range("E38").Value = "  -4.50%  "
$ws.Range("D39").Value = "'0.0813"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("D41").Value = "3.023.92"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").Value = "'8.27"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "'0.269"
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").Value = "'28.19"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'2.26"
$ws.Range("E46").Value = "  +7.54%  "
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "'118.11"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").Value = "0.0₃0513"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "'2.08"
$ws.Range("E51").Value = "  +0.12%  "
